$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 4153
$ws.Range("I80").Value = 305.9375
$ws.Range("J80").Value = 8549.643
$ws.Range("K80").Value = 917.8125
$ws.Range("L80").Value = 25648.929
$ws.Range("M80").Value = 80.1875
$ws.Range("N80").Value = -27644.929
$ws.Range("H83").Value = 4153
$ws.Range("I83").Value = 305.9375
$ws.Range("J83").Value = 8549.643
$ws.Range("K83").Value = 2753.4375
$ws.Range("L83").Value = 76946.787
$ws.Range("M83").Value = 2238.5625
$ws.Range("N83").Value = -86930.787
$ws.Range("H125").Value = 7057.8125
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 7057.8125
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").Value = 63520.3125
$ws.Range("N125").Value = -68440.3125
$ws.Range("H129").Value = 1033
$ws.Range("J129").Value = 1117.1786
$ws.Range("L129").Value = 3351.5358
$ws.Range("N129").Value = -13351.5358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 29763.8
$ws.Range("J80").Value = 39939.668
$ws.Range("L80").Value = 39939.668
$ws.Range("N80").Value = -41935.668
$ws.Range("H83").Value = 29763.8
$ws.Range("J83").Value = 39939.668
$ws.Range("L83").Value = 119819.004
$ws.Range("N83").Value = -129803.004
$ws.Range("H102").Value = 3370.9092
$ws.Range("I102").Value = 2840
$ws.Range("J102").Value = 4300
$ws.Range("K102").Value = 2840
$ws.Range("L102").Value = 4300
$ws.Range("M102").Value = -1218
$ws.Range("N102").Value = -7544
$ws.Range("H122").Value = 1600.5333
$ws.Range("I122").Value = 1310.8
$ws.Range("K122").Value = 3932.4
$ws.Range("M122").Value = -1482.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19226.111
$ws.Range("I82").Value = 2382.6
$ws.Range("J82").Value = 40280.5
$ws.Range("K82").Value = 2382.6
$ws.Range("L82").Value = 40280.5
$ws.Range("M82").Value = -1999.6
$ws.Range("N82").Value = -41046.5
$ws.Range("H85").Value = 19226.111
$ws.Range("I85").Value = 2382.6
$ws.Range("J85").Value = 40280.5
$ws.Range("K85").Value = 2382.6
$ws.Range("L85").Value = 40280.5
$ws.Range("M85").Value = -1056.6
$ws.Range("N85").Value = -42932.5
$ws.Range("H99").Value = 500
$ws.Range("I99").Value = 500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 998
$ws.Range("H134").Value = 3451.611
$ws.Range("I134").Value = 3580.3845
$ws.Range("J134").Value = 3116.8
$ws.Range("K134").Value = 10741.1535
$ws.Range("L134").Value = 9350.400000000001
$ws.Range("M134").Value = -8206.1535
$ws.Range("N134").Value = -14420.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5501.303
$ws.Range("I31").Value = 6588.8945
$ws.Range("J31").Value = 4025.2856
$ws.Range("K31").Value = 6588.8945
$ws.Range("L31").Value = 4025.2856
$ws.Range("M31").Value = -6293.8945
$ws.Range("N31").Value = -4615.2856
$ws.Range("H34").Value = 5501.303
$ws.Range("I34").Value = 6588.8945
$ws.Range("J34").Value = 4025.2856
$ws.Range("K34").Value = 6588.8945
$ws.Range("L34").Value = 4025.2856
$ws.Range("M34").Value = -6386.8945
$ws.Range("N34").Value = -4429.2856
$ws.Range("H122").Value = 8539.895
$ws.Range("I122").Value = 3801.875
$ws.Range("J122").Value = 33809.332
$ws.Range("K122").Value = 11405.625
$ws.Range("L122").Value = 101427.996
$ws.Range("M122").Value = -8955.625
$ws.Range("N122").Value = -106327.996
$ws.Range("H134").Value = 2318.318
$ws.Range("I134").Value = 1363.675
$ws.Range("J134").Value = 3787
$ws.Range("K134").Value = 4091.025
$ws.Range("L134").Value = 11361
$ws.Range("M134").Value = -1556.025
$ws.Range("N134").Value = -16431

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8309.632
$ws.Range("I3").Value = 5049.091
$ws.Range("J3").Value = 12792.875
$ws.Range("K3").Value = 15147.273
$ws.Range("L3").Value = 38378.625
$ws.Range("M3").Value = -15035.273
$ws.Range("N3").Value = -38602.625
$ws.Range("H4").Value = 250.25
$ws.Range("I4").Value = 250.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 750.75
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -638.75
$ws.Range("H11").Value = 226.36363
$ws.Range("I11").Value = 210
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 630
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -490
$ws.Range("N11").Value = -1180
$ws.Range("H18").Value = 40000284
$ws.Range("I18").Value = 66666750
$ws.Range("J18").Value = 583
$ws.Range("K18").Value = 200000250
$ws.Range("L18").Value = 1749
$ws.Range("M18").Value = -200000081
$ws.Range("N18").Value = -2087
$ws.Range("H31").Value = 250
$ws.Range("I31").Value = 250
$ws.Range("K31").Value = 750
$ws.Range("M31").Value = -462
$ws.Range("H107").Value = 895.0833
$ws.Range("I107").Value = 330.6111
$ws.Range("J107").Value = 1459.5555
$ws.Range("K107").Value = 991.8333
$ws.Range("L107").Value = 4378.666499999999
$ws.Range("M107").Value = 928.1667
$ws.Range("N107").Value = -8218.666499999999
$ws.Range("H122").Value = 862.54285
$ws.Range("I122").Value = 335.7143
$ws.Range("J122").Value = 994.25
$ws.Range("K122").Value = 3021.4287
$ws.Range("L122").Value = 8948.25
$ws.Range("M122").Value = -571.4286999999999
$ws.Range("N122").Value = -13848.25
$ws.Range("H127").Value = 3153.4243
$ws.Range("J127").Value = 3153.4243
$ws.Range("L127").Value = 9460.2729
$ws.Range("N127").Value = -19380.2729
$ws.Range("H129").Value = 3385.8572
$ws.Range("I129").Value = 2283.5
$ws.Range("J129").Value = 10000
$ws.Range("K129").Value = 6850.5
$ws.Range("L129").Value = 30000
$ws.Range("M129").Value = -1850.5
$ws.Range("N129").Value = -40000
$ws.Range("H139").Value = 1905920
$ws.Range("I139").Value = 3708898.2
$ws.Range("J139").Value = 2776.1667
$ws.Range("K139").Value = 11126694.6
$ws.Range("L139").Value = 8328.500100000001
$ws.Range("M139").Value = -11121554.6
$ws.Range("N139").Value = -18608.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 366705.78
$ws.Range("I61").Value = 11633.728
$ws.Range("J61").Value = 1668636.6
$ws.Range("K61").Value = 11633.728
$ws.Range("L61").Value = 1668636.6
$ws.Range("M61").Value = -11431.728
$ws.Range("N61").Value = -1669040.6
$ws.Range("H100").Value = 3900.0476
$ws.Range("I100").Value = 3562.5625
$ws.Range("J100").Value = 4980
$ws.Range("K100").Value = 3562.5625
$ws.Range("L100").Value = 4980
$ws.Range("M100").Value = -3021.5625
$ws.Range("N100").Value = -6062
$ws.Range("H113").Value = 366705.78
$ws.Range("I113").Value = 11633.728
$ws.Range("J113").Value = 1668636.6
$ws.Range("K113").Value = 11633.728
$ws.Range("L113").Value = 1668636.6
$ws.Range("M113").Value = -9463.727999999999
$ws.Range("N113").Value = -1672976.6
$ws.Range("H122").Value = 6440.0576
$ws.Range("I122").Value = 5719.3716
$ws.Range("K122").Value = 17158.1148
$ws.Range("M122").Value = -14708.1148

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1237.5238
$ws.Range("I126").Value = 1282.6666
$ws.Range("J126").Value = 966.6667
$ws.Range("K126").Value = 3847.9998
$ws.Range("L126").Value = 2900.0001
$ws.Range("M126").Value = -1377.9998
$ws.Range("N126").Value = -7840.0001
$ws.Range("H132").Value = 1875.9025
$ws.Range("I132").Value = 1104.9565
$ws.Range("J132").Value = 2861
$ws.Range("K132").Value = 3314.8695
$ws.Range("L132").Value = 8583
$ws.Range("M132").Value = -784.8694999999998
$ws.Range("N132").Value = -13643
$ws.Range("H136").Value = 8828.543
$ws.Range("I136").Value = 7414.4736
$ws.Range("J136").Value = 10507.75
$ws.Range("K136").Value = 22243.4208
$ws.Range("L136").Value = 31523.25
$ws.Range("M136").Value = -19693.4208
$ws.Range("N136").Value = -36623.25
